# Replace the locale-dependent TEXT(M,"MM/AAAA") / TEXT(M,"AAAAMM") formulas
# (which rely on the Spanish "AAAA" year token) with locale-independent
# formulas that build the month/year string from the already-computed
# O (YEAR) and P (MONTH) helper columns. This matches the workbook across
# regional settings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 has its own (non-shared) formulas.
$ws.Range("N2").Formula = '=TEXT(P2,"00")&"/"&O2'
$ws.Range("S2").Formula = '=CONCATENATE(B2," - ",O2,TEXT(P2,"00")," - ",IF(F2="",A2,F2))'

# Rows 3:14 share formulas (N3:N14 and S3:S14, anchored on row 3).
$ws.Range("N3:N14").Formula = '=TEXT(P3,"00")&"/"&O3'
$ws.Range("S3:S14").Formula = '=CONCATENATE(B3," - ",O3,TEXT(P3,"00")," - ",IF(F3="",A3,F3))'

# Row 6 had its own standalone copy of the S formula (separate from the
# S3:S14 shared group); refresh it with the same pattern so it stays in
# sync with the rest of the column.
$ws.Range("S6").Formula = '=CONCATENATE(B6," - ",O6,TEXT(P6,"00")," - ",IF(F6="",A6,F6))'
